$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the crypto list rows 2-51
# Values are written as text to preserve original formatting (e.g. "27.872.58", "  -0.25%  ").
# Purely-numeric-looking price values are forced to text with a leading apostrophe and then
# restored to the default (unstyled) cell style so the underlying XML has no numFmt/style change.

$ws.Range('D2').Value = '27.872.58'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '1.633.02'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('D6').Value = "'0.520"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = "'23.36"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('E9').Value = '  -1.14%  '
$ws.Range('D10').Value = "'0.0611"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.25%  '
$ws.Range('D11').Value = "'0.0883"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('D12').Value = '1.864.40'
$ws.Range('E12').Value = '  -0.27%  '
$ws.Range('D13').Value = '1.642.83'
$ws.Range('E13').Value = '  +0.53%  '
$ws.Range('D14').Value = "'4.02"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.18%  '
$ws.Range('D15').Value = "'0.561"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.21%  '
$ws.Range('D16').Value = "'65.32"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').Value = '27.880.66'
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('D18').Value = "'228.73"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.04%  '
$ws.Range('D19').Value = "'7.67"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.83%  '
$ws.Range('E20').Value = '  -0.08%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('E22').Value = '  -0.89%  '
$ws.Range('D23').Value = "'10.01"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.95%  '
$ws.Range('D24').Value = "'2.08"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').Value = "'155.12"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('D26').Value = "'6.87"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.51%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').Value = "'15.53"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.59%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('D31').Value = "'0.0481"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('E33').Value = '  +1.15%  '
$ws.Range('D34').Value = '1.393.02'
$ws.Range('E34').Value = '  -1.07%  '
$ws.Range('E35').Value = '  +0.83%  '
$ws.Range('E36').Value = '  +9.31%  '
$ws.Range('E37').Value = '  -0.66%  '
$ws.Range('D38').Value = "'0.0170"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.56%  '
$ws.Range('E39').Value = '  -0.69%  '
$ws.Range('D40').Value = "'0.848"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.16%  '
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('E42').Value = '  -0.97%  '
$ws.Range('D43').Value = "'65.86"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.94%  '
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('D45').Value = "'5.44"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.48%  '
$ws.Range('D46').Value = '1.773.65'
$ws.Range('E46').Value = '  -0.32%  '
$ws.Range('E47').Value = '  -2.83%  '
$ws.Range('D48').Value = "'88.75"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.85%  '
$ws.Range('E49').Value = '  +2.04%  '
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('D51').Value = "'7.65"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.15%  '
